# Add "2022-Q1" sheet (new fund-holding detail) positioned right before "总计",
# and insert a new first data-row into "总计" summarizing the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet just before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$ns = $wb.Worksheets.Add($totalSheet)
$ns.Name = "2022-Q1"

# Reuse the header / index-column look-and-feel (bold + bordered style) that
# every other quarterly sheet already uses, by copying formats from the
# "2021-Q4" sheet (same visual style carried across all quarter tabs).
$styleSrc = $wb.Worksheets.Item("2021-Q4")

$styleSrc.Range("B1:H1").Copy()
$ns.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$styleSrc.Range("A2").Copy()
$ns.Range("A2:A9").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$ns.Cells.Item(1,2).Value = "基金代码"
$ns.Cells.Item(1,3).Value = "基金名称"
$ns.Cells.Item(1,4).Value = "基金规模"
$ns.Cells.Item(1,5).Value = "股票总仓位"
$ns.Cells.Item(1,6).Value = "仓位占比"
$ns.Cells.Item(1,7).Value = "持有市值(亿元)"
$ns.Cells.Item(1,8).Value = "仓位排名"

# Data rows — columns B..G carry text (fund codes / figures), stored as text
# (same as every other quarterly sheet) so leading zeros in fund codes and
# exact figure formatting survive; column A (index) and H (rank) are numeric.
$data = @(
  @(0, "512710", "富国中证军工龙头ETF",                         "42.62", "99.44", "2.81", "1.1976", 10),
  @(1, "167506", "安信中证深圳科技创新主题指数（LOF）A",         "1.27",  "90.25", "4.01", "0.0509", 7),
  @(2, "012216", "红塔红土盛利混合型证券投资基金A",               "2.82",  "51.68", "0.88", "0.0248", 6),
  @(3, "010756", "兴华永兴混合A",                                 "0.35",  "94.57", "4.04", "0.0141", 8),
  @(4, "167507", "安信中证深圳科技创新主题指数（LOF）C",         "0.33",  "90.25", "4.01", "0.0132", 7),
  @(5, "012217", "红塔红土盛利混合型证券投资基金C",               "0.53",  "51.68", "0.88", "0.0047", 6),
  @(6, "002303", "金鹰智慧生活灵活配置混合",                     "0.11",  "89.88", "3.23", "0.0036", 8),
  @(7, "010757", "兴华永兴混合C",                                 "0.01",  "94.57", "4.04", "0.0004", 8)
)

$r = 2
foreach ($row in $data) {
    $ns.Cells.Item($r, 1).Value = $row[0]

    $textRange = $ns.Range("B$r`:G$r")
    $textRange.NumberFormat = "@"
    $ns.Cells.Item($r, 2).Value = $row[1]
    $ns.Cells.Item($r, 3).Value = $row[2]
    $ns.Cells.Item($r, 4).Value = $row[3]
    $ns.Cells.Item($r, 5).Value = $row[4]
    $ns.Cells.Item($r, 6).Value = $row[5]
    $ns.Cells.Item($r, 7).Value = $row[6]
    $textRange.Style = "Normal"

    $ns.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------------
# 2) Insert the corresponding summary row into "总计", ahead of "2021-Q4"
# ---------------------------------------------------------------------------
$ts = $wb.Worksheets.Item("总计")
$ts.Rows.Item(2).Insert()

$ts.Cells.Item(2,1).Value = 0
$ts.Cells.Item(2,2).Value = "2022-Q1"
$ts.Cells.Item(2,3).Value = 8
$ts.Cells.Item(2,4).Value = 1.31

# The row-insert leaves stray formatting behind; restore it to match the
# other rows: column A keeps the bold/bordered index style, B:D stay plain.
$ts.Cells.Item(3,1).Copy()
$ts.Cells.Item(2,1).PasteSpecial(-4122)  # xlPasteFormats
$ts.Range("B2:D2").Style = "Normal"

# Column A is a plain 0-based row counter, independent of the date label —
# renumber it for every row now that a row was inserted at the top.
for ($i = 0; $i -lt 6; $i++) {
    $ts.Cells.Item($i + 2, 1).Value = $i
}
